$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had columns: A=Name, B=ID, C=Amount.
# A new "City" column is inserted as the new column B, pushing the old
# ID column to C and the old Amount column to D.
$lastRow = $ws.UsedRange.Rows.Count

# Shift the Amount column (C -> D) and the ID column (B -> C).
# Range.Copy is used instead of reading/writing .Value so the long
# numeric-looking ID strings keep their original text type instead of
# being coerced into (and losing precision as) numbers, and so that no
# new cell styles get introduced in the process.
$ws.Range("C1:C$lastRow").Copy($ws.Range("D1:D$lastRow"))
$ws.Range("B1:B$lastRow").Copy($ws.Range("C1:C$lastRow"))

# Fill the new City column (B) with each educator's city.
# Row 5 mirrors row 3's educator (Milica Jakšić / same ID) but, matching
# the source data, keeps the city from row 4 ("Subotica") rather than
# row 3's own city ("Zrenjanin").
$ws.Range("B1").Value = "Beograd"
$ws.Range("B2").Value = "Novi Sad"
$ws.Range("B3").Value = "Zrenjanin"
$ws.Range("B4").Value = "Subotica"
$ws.Range("B5").Value = "Subotica"
